$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = '@'
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '26.271.14'
Set-TextValue 'E2' '  -7.61%  '
Set-TextValue 'D3' '1.674.17'
Set-TextValue 'E3' '  -6.06%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  +0.47%  '
Set-TextValue 'D5' '218.28'
Set-TextValue 'E5' '  -5.00%  '
Set-TextValue 'D6' '0.5092'
Set-TextValue 'E6' '  -12.89%  '
Set-TextValue 'D7' '1.006'
Set-TextValue 'E7' '  +0.37%  '
Set-TextValue 'D8' '0.2651'
Set-TextValue 'E8' '  -3.34%  '
Set-TextValue 'D9' '22.03'
Set-TextValue 'E9' '  -4.77%  '
Set-TextValue 'D10' '0.06320'
Set-TextValue 'E10' '  -5.24%  '
Set-TextValue 'D11' '0.07357'
Set-TextValue 'E11' '  -2.19%  '
Set-TextValue 'D12' '1.680.51'
Set-TextValue 'E12' '  -5.82%  '
Set-TextValue 'D13' '4.534'
Set-TextValue 'E13' '  -4.42%  '
Set-TextValue 'D14' '0.5732'
Set-TextValue 'E14' '  -5.32%  '
Set-TextValue 'D15' '1.908.33'
Set-TextValue 'E15' '  -5.67%  '
Set-TextValue 'D16' '0.000008526'
Set-TextValue 'E16' '  -0.94%  '
Set-TextValue 'D17' '64.65'
Set-TextValue 'E17' '  -13.47%  '
Set-TextValue 'D18' '26.374.16'
Set-TextValue 'E18' '  -7.13%  '
Set-TextValue 'D19' '4.993'
Set-TextValue 'E19' '  -6.92%  '
Set-TextValue 'E20' '  +0.20%  '
Set-TextValue 'D21' '10.84'
Set-TextValue 'E21' '  -4.30%  '
Set-TextValue 'D22' '185.95'
Set-TextValue 'E22' '  -9.91%  '
Set-TextValue 'D23' '6.214'
Set-TextValue 'E23' '  -7.77%  '
Set-TextValue 'D24' '1.006'
Set-TextValue 'E24' '  +0.42%  '
Set-TextValue 'D25' '143.47'
Set-TextValue 'E25' '  -5.54%  '
Set-TextValue 'D26' '7.523'
Set-TextValue 'E26' '  -6.90%  '
Set-TextValue 'D27' '0.1170'
Set-TextValue 'E27' '  -5.97%  '
Set-TextValue 'D28' '15.70'
Set-TextValue 'E28' '  -3.19%  '
Set-TextValue 'D29' '1.332'
Set-TextValue 'E29' '  -5.08%  '
Set-TextValue 'D30' '0.05813'
Set-TextValue 'E30' '  -4.95%  '
Set-TextValue 'D31' '1.328'
Set-TextValue 'E31' '  -5.90%  '
Set-TextValue 'D32' '3.505'
Set-TextValue 'E32' '  -6.59%  '
Set-TextValue 'D33' '3.498'
Set-TextValue 'E33' '  -6.81%  '
Set-TextValue 'D34' '1.647'
Set-TextValue 'E34' '  -0.95%  '
Set-TextValue 'D35' '1.005'
Set-TextValue 'E35' '  -3.47%  '
Set-TextValue 'D36' '0.5937'
Set-TextValue 'E36' '  -6.31%  '
Set-TextValue 'D37' '2.358'
Set-TextValue 'E37' '  -5.74%  '
Set-TextValue 'D38' '2.660'
Set-TextValue 'E38' '  -0.92%  '
Set-TextValue 'B39' 'Maker'
Set-TextValue 'C39' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D39' '1.097.71'
Set-TextValue 'E39' '  -3.61%  '
Set-TextValue 'B40' 'VeChain'
Set-TextValue 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D40' '0.01602'
Set-TextValue 'E40' '  -4.00%  '
Set-TextValue 'D41' '5.890'
Set-TextValue 'E41' '  -6.02%  '
Set-TextValue 'D42' '0.8606'
Set-TextValue 'E42' '  -1.46%  '
Set-TextValue 'D43' '1.005'
Set-TextValue 'E43' '  +0.01%  '
Set-TextValue 'D44' '99.67'
Set-TextValue 'E44' '  -0.62%  '
Set-TextValue 'D45' '1.835.15'
Set-TextValue 'E45' '  -5.17%  '
Set-TextValue 'D46' '0.00000000114'
Set-TextValue 'E46' '  +4.63%  '
Set-TextValue 'D47' '56.17'
Set-TextValue 'E47' '  -5.65%  '
Set-TextValue 'D48' '1.003'
Set-TextValue 'E48' '  +0.66%  '
Set-TextValue 'D49' '8.039'
Set-TextValue 'E49' '  -3.79%  '
Set-TextValue 'D50' '0.4318'
Set-TextValue 'E50' '  -3.19%  '
Set-TextValue 'D51' '0.05205'
Set-TextValue 'E51' '  -3.85%  '
